# atmPy materials.xlsx: add mixing-rule inputs for H2SO4 / (NH4)HSO4
# (molecular weights + ion composition) and unhide the helper columns
# E:H that hold the lookup tables those mixing rules read from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Unhide the (previously hidden) helper columns E:H ---------------
$ws.Columns("E:H").Hidden = $false

# --- New data: molecular weight for H2SO4 (row 4) and (NH4)HSO4 (row 5)
$ws.Range("I4").Value = 98
$ws.Range("I5").Value = 115

# --- New data: ion composition for (NH4)HSO4 (row 5) -----------------
$ws.Range("L5").Value = "ammonium"
$ws.Range("M5").Value = "hydrogen_sulfate"

# --- Update the selection to match the author's last cursor position --
$ws.Range("I6").Select()
